$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings are not
# auto-converted to numbers by Excel's type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.892.78"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "2.664.25"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "592.24"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").Value = "146.93"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -1.22%  "

$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").Value = "5.63"
$ws.Range("E10").Value = "  -0.57%  "

$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").Value = "27.66"
$ws.Range("E13").Value = "  +1.01%  "

$ws.Range("D14").Value = "3.142.33"
$ws.Range("E14").Value = "  +2.52%  "

$ws.Range("D15").Value = "63.749.66"
$ws.Range("E15").Value = "  +1.17%  "

$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").Value = "2.742.24"
$ws.Range("E17").Value = "  +5.49%  "

$ws.Range("D18").Value = "11.33"
$ws.Range("E18").Value = "  -0.56%  "

$ws.Range("D19").Value = "341.86"
$ws.Range("E19").Value = "  -0.88%  "

$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").Value = "6.81"
$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("D23").Value = "68.24"
$ws.Range("E23").Value = "  +1.46%  "

$ws.Range("E24").Value = "  +10.97%  "

$ws.Range("D25").Value = "1.66"
$ws.Range("E25").Value = "  +3.28%  "

$ws.Range("D26").Value = "0.167"
$ws.Range("E26").Value = "  -1.52%  "

$ws.Range("D27").Value = "551.35"
$ws.Range("E27").Value = "  +17.00%  "

$ws.Range("D28").Value = "8.49"
$ws.Range("E28").Value = "  +0.90%  "

$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("D30").Value = "7.92"
$ws.Range("E30").Value = "  +0.65%  "

$ws.Range("E31").Value = "  +2.82%  "

$ws.Range("E32").Value = "  +9.90%  "

$ws.Range("D33").Value = "0.0₃0817"
$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("D34").Value = "175.23"
$ws.Range("E34").Value = "  -0.72%  "

$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").Value = "0.401"
$ws.Range("E36").Value = "  -0.95%  "

$ws.Range("D37").Value = "19.15"
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("E38").Value = "  +1.63%  "

$ws.Range("D39").Value = "1.76"
$ws.Range("E39").Value = "  +2.61%  "

$ws.Range("D40").Value = "172.44"
$ws.Range("E40").Value = "  +8.07%  "

$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").Value = "40.58"
$ws.Range("E42").Value = "  +2.68%  "

$ws.Range("E43").Value = "  -1.20%  "

$ws.Range("D44").Value = "21.61"
$ws.Range("E44").Value = "  +2.07%  "

$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("E46").Value = "  -0.45%  "

$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("D49").Value = "18.75"
$ws.Range("E49").Value = "  +0.43%  "

$ws.Range("E50").Value = "  +1.11%  "

$ws.Range("E51").Value = "  -0.74%  "

# Restore the original (default/Normal) style on column D so no
# lingering text-format style is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"